# Diary update for Yijia Zhang: add 5 new diary entries (rows 21-25),
# matching the formatting of the existing template row 20, then restyle
# the remaining blank filler rows (26-124) to match the "used" blank-row
# look that appears once real content pushes further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Capture the blank-row formats that already exist (on row 21, still
#    untouched at this point) before we overwrite that row with new
#    diary content. These become the formats used by the blank filler
#    rows from 26 downward.
# ---------------------------------------------------------------------
$ws.Range("B21").Copy()
$ws.Range("A26:A124").PasteSpecial(-4122)

$ws.Range("D21").Copy()
$ws.Range("B26:F124").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Seed rows 21-25 with row 20's formatting (date/text column styles,
#    borders, fills, fonts) so we only need to patch a few per-cell
#    differences afterward.
# ---------------------------------------------------------------------
$ws.Range("A20:G20").Copy()
$ws.Range("A21:G25").PasteSpecial(-4122)

# Only row 21's Goal cell (D21) keeps row 20's smaller-italic "Lecture"
# style; rows 22-24's Goal cells use the regular column style, and row
# 25 has no Goal entry at all (blank cell, same look but no text format).
$ws.Range("C21").Copy()
$ws.Range("D22:D24").PasteSpecial(-4122)

$ws.Range("B17").Copy()
$ws.Range("D25").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Row 21 - 2/11: Study for midterm
# ---------------------------------------------------------------------
$ws.Range("A21").Value = 43872
$ws.Range("A21").NumberFormat = "m/d"
$ws.Range("A21").HorizontalAlignment = -4131
$ws.Range("A21").VerticalAlignment = -4160
$ws.Range("A21").WrapText = $true
$ws.Range("B21").Value = "21:00 - 2.12 ??? "
$ws.Range("C21").Value = "Me"
$ws.Range("D21").Value = "Study for midterm"
$ws.Range("E21").Value = "Reviewed the slides, refreshed memory."
$ws.Range("F21").Value = "Got to review the slides more often in the future, as many details/concepts gets easily forgotten over time. It's extra trouble trying to pick them back up. "
$ws.Range("G21").Value = "THE MIDTERM IS BACK UHHHHH"
$ws.Rows.Item(21).RowHeight = 152.6

# ---------------------------------------------------------------------
# 4. Row 22 - 2/13: Survive the Midterm
# ---------------------------------------------------------------------
$ws.Range("A22").Value = 43874
$ws.Range("A22").NumberFormat = "m/d"
$ws.Range("A22").HorizontalAlignment = -4131
$ws.Range("A22").VerticalAlignment = -4160
$ws.Range("A22").WrapText = $true
$ws.Range("B22").Value = "17:00-19:50"
$ws.Range("C22").Value = "Me"
$ws.Range("D22").Value = "Survive the Midterm"
$ws.Range("E22").Value = "Barely finished the midterm. Quickly went through this week's material."
$ws.Range("F22").Value = "Haven't been able to write quickly since I had never written anything in around half a year" + [char]0x2026 + " Feeling sincerely sorry for Kaj for having to try to read my handwriting. Feels okay about the exam, as there were no multiple choice problems. Could not really focus afterward. Will review the slides at home."
$ws.Range("G22").Value = "Not sure" + [char]0x2026 + " Depends on the midterm score."
$ws.Rows.Item(22).RowHeight = 162.85

# ---------------------------------------------------------------------
# 5. Row 23 - 2/19: Review last week's slides
# ---------------------------------------------------------------------
$ws.Range("A23").Value = 43880
$ws.Range("A23").NumberFormat = "m/d"
$ws.Range("A23").HorizontalAlignment = -4131
$ws.Range("A23").VerticalAlignment = -4160
$ws.Range("A23").WrapText = $true
$ws.Range("B23").Value = "20:00-?"
$ws.Range("C23").Value = "Me"
$ws.Range("D23").Value = "Review last week's slides"
$ws.Range("E23").Value = "Went through the slides"
$ws.Range("F23").Value = "Should have focused in class "
$ws.Range("G23").Value = ":("
$ws.Rows.Item(23).RowHeight = 91.15

# ---------------------------------------------------------------------
# 6. Row 24 - 2/20: Lecture - architectural recovery
# ---------------------------------------------------------------------
$ws.Range("A24").Value = 43881
$ws.Range("A24").NumberFormat = "m/d"
$ws.Range("A24").HorizontalAlignment = -4131
$ws.Range("A24").VerticalAlignment = -4160
$ws.Range("A24").WrapText = $true
$ws.Range("B24").Value = "17:00-19:50"
$ws.Range("C24").Value = "Me"
$ws.Range("D24").Value = "Lecture"
$ws.Range("E24").Value = "Learnt about concepts and methods toward Architectural recovery. "
$ws.Range("F24").Value = "Contents today reminds me of the concept " + [char]0x201C + "architectural drift" + [char]0x201D + " that I learnt in 264P at the beginning of the quarter. The guest speakers were informative and fun. Both (of course) agrees on that separation of concerns is vitally important. Also one should try to write short methods, as having short methods shows that one has thought about and can control his structure and style of code."
$ws.Range("G24").Value = "Positive. "
$ws.Rows.Item(24).RowHeight = 204.1

# ---------------------------------------------------------------------
# 7. Row 25 - 2/24: Score was out
# ---------------------------------------------------------------------
$ws.Range("A25").Value = 43885
$ws.Range("A25").NumberFormat = "m/d"
$ws.Range("A25").HorizontalAlignment = -4131
$ws.Range("A25").VerticalAlignment = -4160
$ws.Range("A25").WrapText = $true
$ws.Range("B25").Value = "??"
$ws.Range("C25").Value = "Me"
$ws.Range("D25").Value = ""
$ws.Range("E25").Value = "Score was out"
$ws.Range("F25").Value = "Maybe I didn't study hard enough?  Not sure what costed my 6 points" + [char]0x2026 + " "
$ws.Range("G25").Value = "Depressed"
$ws.Rows.Item(25).RowHeight = 106.2

# ---------------------------------------------------------------------
# 8. Row 26 becomes taller to match the new layout (first filler row
#    right after the new content).
# ---------------------------------------------------------------------
$ws.Rows.Item(26).RowHeight = 104.45
